# Update "想去人数" (want-to-go count) figures for two events that appear
# on both the "展览" sheet and the "全部类型" sheet.
#   F4: 226 -> 227  (南宁·布谷鸟动漫展5th)
#   F5: 3811 -> 3812 (南宁·2024良牙动漫秋季盛典（秋典）)

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 227
    $ws.Range("F5").Value = 3812
}
